$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date from A7 (it shifts down to the new row 12)
$ws.Range("A7").Value = $null

# Add a new row 12 with the date that used to be in A7 (05/02/2018, serial 43222)
# Copy the existing date cell's formatting first, then set the raw serial value
# so it reuses the same date style instead of creating a new number format.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A12").Value = 43222

# Fill in the Group names for each team's presentation slot
$ws.Range("B2").Value = "LevidPynch"
$ws.Range("B3").Value = "SoundBlaster"
$ws.Range("B4").Value = "Team Game Suite"
$ws.Range("B5").Value = "Team DTM"
$ws.Range("B6").Value = "Greenthumbs"
$ws.Range("B8").Value = "Team SNES"
$ws.Range("B9").Value = "Dank Spots"
$ws.Range("B10").Value = "Cyber Fox Games"
$ws.Range("B11").Value = "Speed Daemons"

# Update the active selection to match the edited workbook
$ws.Range("D5").Select()
